$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, pushing the existing rows 56-58 down to 57-59.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly data point
# (date 2023-03-15 / serial 45013, volume 25, prices 15000 / 1500).
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C56").Value = 'Los Lagos'
$ws.Range("D56").Value = 45013
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = 100112012
$ws.Range("G56").Value = 'Espinaca'
$ws.Range("H56").Value = 'Sin especificar'
$ws.Range("I56").Value = 'Primera'
$ws.Range("J56").Value = 25
$ws.Range("K56").Value = 15000
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = 15000
$ws.Range("N56").Value = '$/cuna 10 kilos'
$ws.Range("O56").Value = 'Región Metropolitana'
$ws.Range("P56").Value = 1500
$ws.Range("Q56").Value = 10
$ws.Range("R56").Value = 'Hortaliza'
